$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Regenerated "K" values (strikeouts) replacing old "Strike#" derived values in column G
$kValues = @{
    2  = 4
    3  = 1
    4  = 4
    5  = 1
    6  = 3
    7  = 0
    8  = 2
    9  = 2
    10 = 6
    11 = 1
    12 = 1
    13 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
